$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.08097566666666667
$ws.Cells.Item(2, 8).Value = 0.242927
$ws.Cells.Item(2, 9).Value = 0.005588990034505014
$ws.Cells.Item(2, 10).Value = 0.005588990034505015
$ws.Cells.Item(2, 13).Value = 1.485259333333333
$ws.Cells.Item(2, 14).Value = 4.455778
$ws.Cells.Item(2, 15).Value = 0.3057455162066235
$ws.Cells.Item(2, 16).Value = 0.3057455162066235
$ws.Cells.Item(2, 17).Value = 0.1202698646895556
$ws.Cells.Item(2, 18).Value = 1.082428782206
$ws.Cells.Item(2, 19).Value = 0.00170880864317341
$ws.Cells.Item(2, 20).Value = 0.00170880864317341

$ws.Cells.Item(3, 7).Value = 0.08097566666666667
$ws.Cells.Item(3, 8).Value = 0.242927
$ws.Cells.Item(3, 9).Value = 0.005588990034505014
$ws.Cells.Item(3, 10).Value = 0.005588990034505015
$ws.Cells.Item(3, 15).Value = 0.2805555239151429
$ws.Cells.Item(3, 16).Value = 0.2805555239151429
$ws.Cells.Item(3, 17).Value = 0.1103609803271111
$ws.Cells.Item(3, 18).Value = 0.9932488229440002
$ws.Cells.Item(3, 19).Value = 0.001568022027287067
$ws.Cells.Item(3, 20).Value = 0.001568022027287067

$ws.Cells.Item(4, 7).Value = 0.08097566666666667
$ws.Cells.Item(4, 8).Value = 0.242927
$ws.Cells.Item(4, 9).Value = 0.005588990034505014
$ws.Cells.Item(4, 10).Value = 0.005588990034505015
$ws.Cells.Item(4, 15).Value = 0.4136989598782336
$ws.Cells.Item(4, 16).Value = 0.4136989598782336
$ws.Cells.Item(4, 17).Value = 0.1627350698191111
$ws.Cells.Item(4, 18).Value = 1.464615628372
$ws.Cells.Item(4, 19).Value = 0.002312159364044537
$ws.Cells.Item(4, 20).Value = 0.002312159364044537

$ws.Cells.Item(5, 9).Value = 0.6976944377922635
$ws.Cells.Item(5, 10).Value = 0.6976944377922635
$ws.Cells.Item(5, 13).Value = 1.485259333333333
$ws.Cells.Item(5, 14).Value = 4.455778
$ws.Cells.Item(5, 15).Value = 0.3057455162066235
$ws.Cells.Item(5, 16).Value = 0.3057455162066235
$ws.Cells.Item(5, 17).Value = 15.013735059444
$ws.Cells.Item(5, 18).Value = 135.123615534996
$ws.Cells.Item(5, 19).Value = 0.2133169460372856
$ws.Cells.Item(5, 20).Value = 0.2133169460372856

$ws.Cells.Item(6, 9).Value = 0.6976944377922635
$ws.Cells.Item(6, 10).Value = 0.6976944377922635
$ws.Cells.Item(6, 15).Value = 0.2805555239151429
$ws.Cells.Item(6, 16).Value = 0.2805555239151429
$ws.Cells.Item(6, 19).Value = 0.1957420285274896
$ws.Cells.Item(6, 20).Value = 0.1957420285274896

$ws.Cells.Item(7, 9).Value = 0.6976944377922635
$ws.Cells.Item(7, 10).Value = 0.6976944377922635
$ws.Cells.Item(7, 15).Value = 0.4136989598782336
$ws.Cells.Item(7, 16).Value = 0.4136989598782336
$ws.Cells.Item(7, 19).Value = 0.2886354632274884
$ws.Cells.Item(7, 20).Value = 0.2886354632274884

$ws.Cells.Item(8, 9).Value = 0.2967165721732315
$ws.Cells.Item(8, 10).Value = 0.2967165721732316
$ws.Cells.Item(8, 13).Value = 1.485259333333333
$ws.Cells.Item(8, 14).Value = 4.455778
$ws.Cells.Item(8, 15).Value = 0.3057455162066235
$ws.Cells.Item(8, 16).Value = 0.3057455162066235
$ws.Cells.Item(8, 17).Value = 6.385064522589334
$ws.Cells.Item(8, 18).Value = 57.46558070330401
$ws.Cells.Item(8, 19).Value = 0.09071976152616452
$ws.Cells.Item(8, 20).Value = 0.09071976152616454

$ws.Cells.Item(9, 9).Value = 0.2967165721732315
$ws.Cells.Item(9, 10).Value = 0.2967165721732316
$ws.Cells.Item(9, 15).Value = 0.2805555239151429
$ws.Cells.Item(9, 16).Value = 0.2805555239151429
$ws.Cells.Item(9, 19).Value = 0.08324547336036628
$ws.Cells.Item(9, 20).Value = 0.08324547336036629

$ws.Cells.Item(10, 9).Value = 0.2967165721732315
$ws.Cells.Item(10, 10).Value = 0.2967165721732316
$ws.Cells.Item(10, 15).Value = 0.4136989598782336
$ws.Cells.Item(10, 16).Value = 0.4136989598782336
$ws.Cells.Item(10, 19).Value = 0.1227513372867007
$ws.Cells.Item(10, 20).Value = 0.1227513372867007

Write-Output "Applied TPM update to LR-pairs sheet"
